$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C9 value from "Y" to "Y (OK)"
$ws.Range("C9").Value = "Y (OK)"

# Change the active cell selection to C11
$ws.Range("C11").Select()

# Resize the workbook window
$excel.ActiveWindow.Width = 14000
